# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 716
$ws.Range("F17").Value = 1057
$ws.Range("F20").Value = 420
$ws.Range("F23").Value = 674
$ws.Range("F25").Value = 4150
$ws.Range("F27").Value = 2671

# --- 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F19").Value = 20
$ws.Range("F20").Value = 20
$ws.Range("F50").Value = 313

# --- 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 148
$ws.Range("F14").Value = 763

# --- 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 148
$ws.Range("F8").Value = 763
$ws.Range("F22").Value = 1057
$ws.Range("F25").Value = 20
$ws.Range("F30").Value = 420
$ws.Range("F38").Value = 674
$ws.Range("F43").Value = 4150
$ws.Range("F45").Value = 2671
